# Auto-generated PowerShell COM-interop script
# Applies numeric updates to columns H-N across multiple sheets
# as described by the target OOXML diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4632179.5
$ws.Range("I18").Value = 9260759
$ws.Range("J18").Value = 3600
$ws.Range("K18").Value = 9260759
$ws.Range("L18").Value = 3600
$ws.Range("M18").Value = -9260475
$ws.Range("N18").Value = -4168

$ws.Range("H113").Value = 7324.3335
$ws.Range("I113").Value = 4980
$ws.Range("J113").Value = 9668.666999999999
$ws.Range("K113").Value = 4980
$ws.Range("L113").Value = 9668.666999999999
$ws.Range("M113").Value = -1726
$ws.Range("N113").Value = -16176.667

$ws.Range("H129").Value = 5103134
$ws.Range("I129").Value = 83335540
$ws.Range("J129").Value = 1021.2174
$ws.Range("K129").Value = 250006620
$ws.Range("L129").Value = 3063.6522
$ws.Range("M129").Value = -250001620
$ws.Range("N129").Value = -13063.6522


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = 0

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = 0

$ws.Range("H80").Value = 29405.6
$ws.Range("J80").Value = 29405.6
$ws.Range("L80").Value = 29405.6
$ws.Range("N80").Value = -31401.6

$ws.Range("H83").Value = 29405.6
$ws.Range("J83").Value = 29405.6
$ws.Range("L83").Value = 88216.79999999999
$ws.Range("N83").Value = -98200.79999999999

$ws.Range("H122").Value = 3085.697
$ws.Range("I122").Value = 2121.12
$ws.Range("J122").Value = 6100
$ws.Range("K122").Value = 6363.36
$ws.Range("L122").Value = 18300
$ws.Range("M122").Value = -3913.36
$ws.Range("N122").Value = -23200

$ws.Range("H135").Value = 19496.943
$ws.Range("J135").Value = 19496.943
$ws.Range("L135").Value = 19496.943
$ws.Range("N135").Value = -29636.943

$ws.Range("H139").Value = 16624.162
$ws.Range("J139").Value = 16624.162
$ws.Range("L139").Value = 16624.162
$ws.Range("N139").Value = -26904.162


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H82").Value = 20159.646
$ws.Range("I82").Value = 4741.5
$ws.Range("J82").Value = 28569.545
$ws.Range("K82").Value = 4741.5
$ws.Range("L82").Value = 28569.545
$ws.Range("M82").Value = -4358.5
$ws.Range("N82").Value = -29335.545

$ws.Range("H85").Value = 20159.646
$ws.Range("I85").Value = 4741.5
$ws.Range("J85").Value = 28569.545
$ws.Range("K85").Value = 4741.5
$ws.Range("L85").Value = 28569.545
$ws.Range("M85").Value = -3415.5
$ws.Range("N85").Value = -31221.545

$ws.Range("H96").Value = 29000
$ws.Range("I96").Value = 24000
$ws.Range("J96").Value = 34000
$ws.Range("K96").Value = 24000
$ws.Range("L96").Value = 34000
$ws.Range("M96").Value = -21254
$ws.Range("N96").Value = -39492

$ws.Range("H132").Value = 21286.285
$ws.Range("J132").Value = 21286.285
$ws.Range("L132").Value = 21286.285
$ws.Range("N132").Value = -31406.285

$ws.Range("H135").Value = 25875.268
$ws.Range("J135").Value = 25818.455
$ws.Range("L135").Value = 25818.455
$ws.Range("N135").Value = -35958.455


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30472

$ws.Range("H30").Value = 30000
$ws.Range("J30").Value = 30000
$ws.Range("L30").Value = 30000
$ws.Range("N30").Value = -30182

$ws.Range("H31").Value = 1590313.1
$ws.Range("I31").Value = 2224098.8
$ws.Range("J31").Value = 5848.5557
$ws.Range("K31").Value = 2224098.8
$ws.Range("L31").Value = 5848.5557
$ws.Range("M31").Value = -2223803.8
$ws.Range("N31").Value = -6438.5557

$ws.Range("H34").Value = 1590313.1
$ws.Range("I34").Value = 2224098.8
$ws.Range("J34").Value = 5848.5557
$ws.Range("K34").Value = 2224098.8
$ws.Range("L34").Value = 5848.5557
$ws.Range("M34").Value = -2223896.8
$ws.Range("N34").Value = -6252.5557

$ws.Range("H99").Value = 7451.625
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H122").Value = 4357.091
$ws.Range("I122").Value = 3615.5
$ws.Range("J122").Value = 6334.6665
$ws.Range("K122").Value = 10846.5
$ws.Range("L122").Value = 19003.9995
$ws.Range("M122").Value = -8396.5
$ws.Range("N122").Value = -23903.9995

$ws.Range("H126").Value = 7451.625
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

$ws.Range("H135").Value = 23023.166
$ws.Range("J135").Value = 23023.166
$ws.Range("L135").Value = 23023.166
$ws.Range("N135").Value = -33163.166

$ws.Range("H138").Value = 21722.316
$ws.Range("J138").Value = 21722.316
$ws.Range("L138").Value = 21722.316
$ws.Range("N138").Value = -32002.316

$ws.Range("H140").Value = 20626.703
$ws.Range("J140").Value = 20626.703
$ws.Range("L140").Value = 20626.703
$ws.Range("N140").Value = -30986.703


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2510.7693
$ws.Range("J39").Value = 2678.3333
$ws.Range("L39").Value = 8034.999899999999
$ws.Range("N39").Value = -8622.999899999999

$ws.Range("H94").Value = 3746.6667
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3871.4285
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 11614.2855
$ws.Range("M94").Value = -5324
$ws.Range("N94").Value = -12966.2855

$ws.Range("H122").Value = 11190.8
$ws.Range("I122").Value = 650.6667
$ws.Range("K122").Value = 5856.0003
$ws.Range("M122").Value = -3406.0003


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10476.857
$ws.Range("J46").Value = 10723
$ws.Range("L46").Value = 10723
$ws.Range("N46").Value = -11035

$ws.Range("H57").Value = 17780.5

$ws.Range("H80").Value = 3264.1667
$ws.Range("I80").Value = 3710
$ws.Range("J80").Value = 2640
$ws.Range("K80").Value = 3710
$ws.Range("L80").Value = 2640
$ws.Range("M80").Value = -2712
$ws.Range("N80").Value = -4636

$ws.Range("H83").Value = 3264.1667
$ws.Range("I83").Value = 3710
$ws.Range("J83").Value = 2640
$ws.Range("K83").Value = 18550
$ws.Range("L83").Value = 13200
$ws.Range("M83").Value = -13558
$ws.Range("N83").Value = -23184

$ws.Range("H122").Value = 10556.444
$ws.Range("I122").Value = 12500
$ws.Range("J122").Value = 10001.143
$ws.Range("K122").Value = 37500
$ws.Range("L122").Value = 30003.429
$ws.Range("M122").Value = -35050
$ws.Range("N122").Value = -34903.429


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1754.4445
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 2127.1428
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 2127.1428
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -2717.1428

$ws.Range("H27").Value = 1754.4445
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 2127.1428
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 2127.1428
$ws.Range("M27").Value = -343
$ws.Range("N27").Value = -2341.1428

$ws.Range("H40").Value = 8747.75
$ws.Range("I40").Value = 9796.6
$ws.Range("J40").Value = 6999.6665
$ws.Range("K40").Value = 9796.6
$ws.Range("L40").Value = 6999.6665
$ws.Range("M40").Value = -9660.6
$ws.Range("N40").Value = -7271.6665


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 23534
$ws.Range("J101").Value = 23534
$ws.Range("L101").Value = 23534
$ws.Range("N101").Value = -30024

$ws.Range("H135").Value = 19646.459
$ws.Range("J135").Value = 19646.459
$ws.Range("L135").Value = 19646.459
$ws.Range("N135").Value = -29786.459

